$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update A2 value from "frt38355" to "TOTAL"
$ws.Range("A2").Value = "TOTAL"

# Remove row 3 (which contained "frt38868")
$ws.Range("A3").EntireRow.Delete() | Out-Null

# Update the selection to A3:A13 with active cell A3
$ws.Range("A3:A13").Select() | Out-Null
